# Add a new "pl_insol" / "Insolation Flux [Earth Flux]" row to the
# column_descriptions worksheet, inserted just above the existing
# pl_eqt row (row 16), pushing the remaining rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 16 (shifts rows 16-27 down to 17-28).
$ws.Rows("16:16").Insert()

# Populate the new row with the new column description entry.
$ws.Range("A16").Value = "pl_insol"
$ws.Range("B16").Value = "Insolation Flux [Earth Flux]"

# Match the author's resulting selection state: the newly inserted row
# is left selected (whole-row selection) after the insert.
$ws.Rows("16:16").Select()
